# LOG IN BUTTON ACTIVATED
# A new registrant ("AGNIVA BHATTACHARJEE") submitted the sign-up form again,
# appending a new row (47) to the registrations sheet. Also, the stray
# "=08420880979" formula that had been typed into the Phone column of the
# previous row (46) is cleared back down to a plain text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Row 46: drop the accidental formula in the Phone column (F), -----
#        keeping the same text shown to the user.
$ws.Range("F46").NumberFormat = "@"
$ws.Range("F46").Value = "08420880979"

# --- 2. Append the new registration as row 47. ---------------------------
$ws.Range("A47").Value = "BSS/676f1e6bc5a4"
$ws.Range("B47").Value = "AGNIVA"
$ws.Range("C47").Value = "BHATTACHARJEE"
$ws.Range("D47").Value = "nei"

# "Pass Out Year" looks numeric -- force text so the leading formatting
# matches the rest of the column (stored as text, not a number).
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "2020"

# Phone column keeps the same odd formula the earlier row had.
$ws.Range("F47").Formula = "=08420880979"

$ws.Range("G47").Value = "bhattacharjee.agniva.jobs@gmail.com"

# "Date of Birth" looks like a date -- force text so it is stored verbatim
# instead of being converted to a date serial number.
$ws.Range("H47").NumberFormat = "@"
$ws.Range("H47").Value = "2202-01-21"

$ws.Range("I47").Value = "IT"
$ws.Range("J47").Value = "Google"
$ws.Range("K47").Value = ""
$ws.Range("L47").Value = ""
$ws.Range("M47").Value = ""
$ws.Range("N47").Value = "fvg67684yh"
